{"js": "// Replace the date line and the two-digit multiplication problems in the\n// table with their new values, preserving all existing run/paragraph\n// formatting (we replace only the text of each matched run).\nconst replacements = [\n  [\"2026-02-17 Tuesday\", \"2026-02-18 Wednesday\"],\n  [\"43\u00d785=\", \"38\u00d727=\"],\n  [\"21\u00d794=\", \"70\u00d784=\"],\n  [\"56\u00d741=\", \"31\u00d720=\"],\n  [\"34\u00d758=\", \"91\u00d736=\"],\n  [\"37\u00d745=\", \"60\u00d763=\"],\n  [\"34\u00d753=\", \"15\u00d766=\"],\n  [\"32\u00d783=\", \"47\u00d748=\"],\n  [\"28\u00d793=\", \"94\u00d776=\"],\n  [\"63\u00d785=\", \"24\u00d782=\"],\n  [\"17\u00d715=\", \"61\u00d749=\"],\n  [\"73\u00d774=\", \"66\u00d786=\"],\n  [\"59\u00d769=\", \"97\u00d768=\"],\n  [\"81\u00d793=\", \"13\u00d713=\"],\n  [\"37\u00d711=\", \"62\u00d727=\"],\n  [\"55\u00d738=\", \"97\u00d764=\"],\n  [\"90\u00d742=\", \"21\u00d787=\"],\n  [\"37\u00d744=\", \"14\u00d764=\"],\n  [\"86\u00d785=\", \"32\u00d760=\"],\n  [\"50\u00d775=\", \"66\u00d775=\"],\n  [\"62\u00d783=\", \"28\u00d725=\"],\n  [\"25\u00d749=\", \"77\u00d779=\"],\n  [\"83\u00d715=\", \"53\u00d731=\"],\n  [\"14\u00d727=\", \"72\u00d760=\"],\n  [\"50\u00d773=\", \"98\u00d781=\"],\n  [\"99\u00d769=\", \"12\u00d795=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and the two-digit multiplication problems in the\n# table with their new values, preserving all existing run/paragraph\n# formatting (Find/Replace only touches the matched text run).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2026-02-17 Tuesday\", \"2026-02-18 Wednesday\"),\n    @(\"43\u00d785=\", \"38\u00d727=\"),\n    @(\"21\u00d794=\", \"70\u00d784=\"),\n    @(\"56\u00d741=\", \"31\u00d720=\"),\n    @(\"34\u00d758=\", \"91\u00d736=\"),\n    @(\"37\u00d745=\", \"60\u00d763=\"),\n    @(\"34\u00d753=\", \"15\u00d766=\"),\n    @(\"32\u00d783=\", \"47\u00d748=\"),\n    @(\"28\u00d793=\", \"94\u00d776=\"),\n    @(\"63\u00d785=\", \"24\u00d782=\"),\n    @(\"17\u00d715=\", \"61\u00d749=\"),\n    @(\"73\u00d774=\", \"66\u00d786=\"),\n    @(\"59\u00d769=\", \"97\u00d768=\"),\n    @(\"81\u00d793=\", \"13\u00d713=\"),\n    @(\"37\u00d711=\", \"62\u00d727=\"),\n    @(\"55\u00d738=\", \"97\u00d764=\"),\n    @(\"90\u00d742=\", \"21\u00d787=\"),\n    @(\"37\u00d744=\", \"14\u00d764=\"),\n    @(\"86\u00d785=\", \"32\u00d760=\"),\n    @(\"50\u00d775=\", \"66\u00d775=\"),\n    @(\"62\u00d783=\", \"28\u00d725=\"),\n    @(\"25\u00d749=\", \"77\u00d779=\"),\n    @(\"83\u00d715=\", \"53\u00d731=\"),\n    @(\"14\u00d727=\", \"72\u00d760=\"),\n    @(\"50\u00d773=\", \"98\u00d781=\"),\n    @(\"99\u00d769=\", \"12\u00d795=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n}\n"}
